$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the "Conversión del día" text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 5.15 = 20427.84 pesos`n✅ 20427.84 pesos = 5.14 = 967.42 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 194
$ws2.Range("O10").Value = 3963
$ws2.Range("N12").Value = 3973.98
$ws2.Range("O12").Value = 188.2
